$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new daily UF / IVP rows (744-774), matching the data-source refresh ---
# First, stamp the row-744:774 block with the same formatting as the last existing
# data row (743) so the new rows share styles (date format in A, number format in B/C)
# without Excel minting brand-new style records.
$ws.Range("A743:C743").Copy()
$ws.Range("A744:C774").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(744, 1).Value = 44206
$ws.Cells.Item(744, 2).Value = 29064.7
$ws.Cells.Item(744, 3).Value = 30252.14
$ws.Cells.Item(745, 1).Value = 44207
$ws.Cells.Item(745, 2).Value = 29067.51
$ws.Cells.Item(745, 3).Value = 30255.02
$ws.Cells.Item(746, 1).Value = 44208
$ws.Cells.Item(746, 2).Value = 29070.32
$ws.Cells.Item(746, 3).Value = 30257.89
$ws.Cells.Item(747, 1).Value = 44209
$ws.Cells.Item(747, 2).Value = 29073.13
$ws.Cells.Item(747, 3).Value = 30260.77
$ws.Cells.Item(748, 1).Value = 44210
$ws.Cells.Item(748, 2).Value = 29075.93
$ws.Cells.Item(748, 3).Value = 30263.64
$ws.Cells.Item(749, 1).Value = 44211
$ws.Cells.Item(749, 2).Value = 29078.74
$ws.Cells.Item(749, 3).Value = 30266.52
$ws.Cells.Item(750, 1).Value = 44212
$ws.Cells.Item(750, 2).Value = 29081.55
$ws.Cells.Item(750, 3).Value = 30269.4
$ws.Cells.Item(751, 1).Value = 44213
$ws.Cells.Item(751, 2).Value = 29084.36
$ws.Cells.Item(751, 3).Value = 30272.27
$ws.Cells.Item(752, 1).Value = 44214
$ws.Cells.Item(752, 2).Value = 29087.18
$ws.Cells.Item(752, 3).Value = 30275.15
$ws.Cells.Item(753, 1).Value = 44215
$ws.Cells.Item(753, 2).Value = 29089.99
$ws.Cells.Item(753, 3).Value = 30278.03
$ws.Cells.Item(754, 1).Value = 44216
$ws.Cells.Item(754, 2).Value = 29092.799999999999
$ws.Cells.Item(754, 3).Value = 30280.9
$ws.Cells.Item(755, 1).Value = 44217
$ws.Cells.Item(755, 2).Value = 29095.61
$ws.Cells.Item(755, 3).Value = 30283.78
$ws.Cells.Item(756, 1).Value = 44218
$ws.Cells.Item(756, 2).Value = 29098.42
$ws.Cells.Item(756, 3).Value = 30286.66
$ws.Cells.Item(757, 1).Value = 44219
$ws.Cells.Item(757, 2).Value = 29101.23
$ws.Cells.Item(757, 3).Value = 30289.54
$ws.Cells.Item(758, 1).Value = 44220
$ws.Cells.Item(758, 2).Value = 29104.04
$ws.Cells.Item(758, 3).Value = 30292.42
$ws.Cells.Item(759, 1).Value = 44221
$ws.Cells.Item(759, 2).Value = 29106.86
$ws.Cells.Item(759, 3).Value = 30295.29
$ws.Cells.Item(760, 1).Value = 44222
$ws.Cells.Item(760, 2).Value = 29109.67
$ws.Cells.Item(760, 3).Value = 30298.17
$ws.Cells.Item(761, 1).Value = 44223
$ws.Cells.Item(761, 2).Value = 29112.48
$ws.Cells.Item(761, 3).Value = 30301.05
$ws.Cells.Item(762, 1).Value = 44224
$ws.Cells.Item(762, 2).Value = 29115.3
$ws.Cells.Item(762, 3).Value = 30303.93
$ws.Cells.Item(763, 1).Value = 44225
$ws.Cells.Item(763, 2).Value = 29118.11
$ws.Cells.Item(763, 3).Value = 30306.81
$ws.Cells.Item(764, 1).Value = 44226
$ws.Cells.Item(764, 2).Value = 29120.92
$ws.Cells.Item(764, 3).Value = 30309.69
$ws.Cells.Item(765, 1).Value = 44227
$ws.Cells.Item(765, 2).Value = 29123.74
$ws.Cells.Item(765, 3).Value = 30312.57
$ws.Cells.Item(766, 1).Value = 44228
$ws.Cells.Item(766, 2).Value = 29126.55
$ws.Cells.Item(766, 3).Value = 30315.45
$ws.Cells.Item(767, 1).Value = 44229
$ws.Cells.Item(767, 2).Value = 29129.37
$ws.Cells.Item(767, 3).Value = 30318.33
$ws.Cells.Item(768, 1).Value = 44230
$ws.Cells.Item(768, 2).Value = 29132.18
$ws.Cells.Item(768, 3).Value = 30321.21
$ws.Cells.Item(769, 1).Value = 44231
$ws.Cells.Item(769, 2).Value = 29135
$ws.Cells.Item(769, 3).Value = 30324.09
$ws.Cells.Item(770, 1).Value = 44232
$ws.Cells.Item(770, 2).Value = 29137.81
$ws.Cells.Item(770, 3).Value = 30326.98
$ws.Cells.Item(771, 1).Value = 44233
$ws.Cells.Item(771, 2).Value = 29140.63
$ws.Cells.Item(771, 3).Value = 30329.86
$ws.Cells.Item(772, 1).Value = 44234
$ws.Cells.Item(772, 2).Value = 29143.439999999999
$ws.Cells.Item(772, 3).Value = 30332.74
$ws.Cells.Item(773, 1).Value = 44235
$ws.Cells.Item(773, 2).Value = 29146.26
$ws.Cells.Item(773, 3).Value = 30335.62
$ws.Cells.Item(774, 1).Value = 44236
$ws.Cells.Item(774, 2).Value = 29149.08
$ws.Cells.Item(774, 3).Value = 30338.5

# --- Column width adjustments for B and C (author widened them in this revision) ---
$ws.Columns.Item(2).ColumnWidth = 19
$ws.Columns.Item(3).ColumnWidth = 24

# --- Row 2 (header row) height reverts to the default (was an explicit 51pt) ---
$ws.Rows.Item(2).AutoFit()

# --- Keep the named range in sync with the newly-extended table ---
$name = $ws.Names.Item(1)
$name.RefersTo = '=UF_IVP_DIARIO!$A$1:$C$774'

# --- Update the view selection to the new last cell, same as a live refresh would ---
[void]$ws.Range("C774").Select()
